$d = $word.ActiveDocument

$replacements = @(
    @("866÷5=", "516÷3="),
    @("663÷3=", "625÷8="),
    @("868÷5=", "647÷8="),
    @("741÷2=", "887÷9="),
    @("993÷5=", "196÷7="),
    @("714÷6=", "332÷6="),
    @("205÷7=", "802÷6="),
    @("301÷8=", "721÷5="),
    @("547÷6=", "778÷3="),
    @("763÷4=", "497÷5="),
    @("656÷6=", "711÷6="),
    @("745÷7=", "255÷6="),
    @("126÷7=", "325÷4="),
    @("230÷3=", "569÷5="),
    @("743÷5=", "920÷2="),
    @("168÷2=", "378÷2="),
    @("177÷4=", "873÷4="),
    @("388÷7=", "420÷2="),
    @("831÷9=", "696÷3="),
    @("523÷5=", "738÷3="),
    @("291÷6=", "227÷5="),
    @("341÷7=", "384÷3="),
    @("406÷9=", "950÷8="),
    @("394÷5=", "908÷4="),
    @("552÷9=", "736÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
